$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: WrappedEther -> Chainlink (coin swapped position, updated data)
$cell = $ws.Range("B15")
$cell.Value = "'" + 'Chainlink'
$cell.Style = "Normal"

$cell = $ws.Range("C15")
$cell.Value = "'" + 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.Value = "'" + '7.335'
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.Value = "'" + '  -1.05%  '
$cell.Style = "Normal"

# Row 16: Chainlink -> WrappedEther
$cell = $ws.Range("B16")
$cell.Value = "'" + 'WrappedEther'
$cell.Style = "Normal"

$cell = $ws.Range("C16")
$cell.Value = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.Value = "'" + '1.799.30'
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.Value = "'" + '  -1.03%  '
$cell.Style = "Normal"

# Row 2
$cell = $ws.Range("D2")
$cell.Value = "'" + '28.218.93'
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.Value = "'" + '  -0.53%  '
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.Value = "'" + '1.807.85'
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.Value = "'" + '  -0.61%  '
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.Value = "'" + '313.15'
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.Value = "'" + '  -0.64%  '
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.Value = "'" + '1.000'
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.Value = "'" + '  -0.32%  '
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.Value = "'" + '0.5142'
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.Value = "'" + '  -1.52%  '
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("D8")
$cell.Value = "'" + '0.3962'
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.Value = "'" + '  +2.77%  '
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("D9")
$cell.Value = "'" + '0.07817'
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.Value = "'" + '  -3.33%  '
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.Value = "'" + '1.109'
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.Value = "'" + '  -0.46%  '
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("D11")
$cell.Value = "'" + '40.82'
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.Value = "'" + '  -2.55%  '
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("D12")
$cell.Value = "'" + '6.378'
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.Value = "'" + '  -0.45%  '
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.Value = "'" + '1.000'
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.Value = "'" + '  -0.32%  '
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.Value = "'" + '20.44'
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.Value = "'" + '  -2.15%  '
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.Value = "'" + '92.97'
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.Value = "'" + '  -1.28%  '
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.Value = "'" + '0.00001079'
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.Value = "'" + '  -2.04%  '
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.Value = "'" + '0.06584'
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.Value = "'" + '  -0.73%  '
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("E20")
$cell.Value = "'" + '  -0.28%  '
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("D21")
$cell.Value = "'" + '17.31'
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.Value = "'" + '  -1.98%  '
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.Value = "'" + '6.019'
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.Value = "'" + '  -0.04%  '
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("D23")
$cell.Value = "'" + '28.257.38'
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.Value = "'" + '  -0.53%  '
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("D24")
$cell.Value = "'" + '11.14'
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.Value = "'" + '  -2.17%  '
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("D25")
$cell.Value = "'" + '2.212'
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.Value = "'" + '  -1.37%  '
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.Value = "'" + '160.64'
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.Value = "'" + '  +0.97%  '
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.Value = "'" + '2.467'
$cell.Style = "Normal"

$cell = $ws.Range("E27")
$cell.Value = "'" + '  +2.56%  '
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("D28")
$cell.Value = "'" + '20.52'
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("D29")
$cell.Value = "'" + '2.013.03'
$cell.Style = "Normal"

$cell = $ws.Range("E29")
$cell.Value = "'" + '  -0.65%  '
$cell.Style = "Normal"

# Row 30
$cell = $ws.Range("D30")
$cell.Value = "'" + '128.31'
$cell.Style = "Normal"

$cell = $ws.Range("E30")
$cell.Value = "'" + '  +3.14%  '
$cell.Style = "Normal"

# Row 31
$cell = $ws.Range("D31")
$cell.Value = "'" + '0.1095'
$cell.Style = "Normal"

$cell = $ws.Range("E31")
$cell.Value = "'" + '  -0.63%  '
$cell.Style = "Normal"

# Row 32
$cell = $ws.Range("D32")
$cell.Value = "'" + '1.060'
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.Value = "'" + '  -1.88%  '
$cell.Style = "Normal"

# Row 33
$cell = $ws.Range("D33")
$cell.Value = "'" + '3.657'
$cell.Style = "Normal"

$cell = $ws.Range("E33")
$cell.Value = "'" + '  -0.55%  '
$cell.Style = "Normal"

# Row 34
$cell = $ws.Range("D34")
$cell.Value = "'" + '5.575'
$cell.Style = "Normal"

$cell = $ws.Range("E34")
$cell.Value = "'" + '  -1.78%  '
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("D35")
$cell.Value = "'" + '0.07160'
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.Value = "'" + '  -3.07%  '
$cell.Style = "Normal"

# Row 36
$cell = $ws.Range("D36")
$cell.Value = "'" + '9.187'
$cell.Style = "Normal"

$cell = $ws.Range("E36")
$cell.Value = "'" + '  +5.33%  '
$cell.Style = "Normal"

# Row 37
$cell = $ws.Range("E37")
$cell.Value = "'" + '  +0.33%  '
$cell.Style = "Normal"

# Row 38
$cell = $ws.Range("D38")
$cell.Value = "'" + '0.2175'
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.Value = "'" + '  -0.94%  '
$cell.Style = "Normal"

# Row 39
$cell = $ws.Range("D39")
$cell.Value = "'" + '5.040'
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.Value = "'" + '  -2.10%  '
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.Value = "'" + '11.54'
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.Value = "'" + '  -6.39%  '
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("D41")
$cell.Value = "'" + '0.6169'
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.Value = "'" + '  -2.41%  '
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("D42")
$cell.Value = "'" + '1.000'
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.Value = "'" + '  -0.34%  '
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.Value = "'" + '1.153'
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.Value = "'" + '  -2.39%  '
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.Value = "'" + '13.15'
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.Value = "'" + '  -1.82%  '
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("D45")
$cell.Value = "'" + '0.5973'
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.Value = "'" + '  -2.44%  '
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("D46")
$cell.Value = "'" + '1.304'
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.Value = "'" + '  -5.56%  '
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.Value = "'" + '3.733'
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.Value = "'" + '  -1.37%  '
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.Value = "'" + '125.06'
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.Value = "'" + '  -1.60%  '
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.Value = "'" + '1.213'
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.Value = "'" + '  +1.23%  '
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.Value = "'" + '1.919'
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.Value = "'" + '  -3.16%  '
$cell.Style = "Normal"

# Row 51
$cell = $ws.Range("D51")
$cell.Value = "'" + '0.06799'
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.Value = "'" + '  -1.34%  '
$cell.Style = "Normal"
